# Replace graph source data on Sheet1 with a "Total Members" header row plus a
# computed Percentage column (D), per "replace graph on my about for full
# html/css version close #1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert two blank rows above the existing table -----------
# (old rows 2-17 become rows 4-19; the black "background" fill in columns
# E:M rides along automatically because Insert() shifts formatted rows)
$ws.Rows("2:3").Insert()

# --- 2. New summary row 2: Total Members / 19 --------------------------------
$ws.Range("B2").Value = "Total Members"
$ws.Range("C2").Value = 19

# --- 3. New header row (row 4): Category | Members | Percentage -------------
$ws.Range("D4").Value = "Percentage"

# --- 4. New Percentage column (D), rows 5-10 ---------------------------------
# D5 alone, then D6:D10 as one fill so Excel records a shared formula (si=0),
# matching the original editor's fill-down.
$ws.Range("D5").Formula = "=(C5/`$C`$12)"
$ws.Range("D6:D10").Formula = "=(C6/`$C`$12)"

# --- 5. Totals row (row 12): sums of Members and Percentage -----------------
$ws.Range("C12").Formula = "=SUM(C5:C10)"
$ws.Range("D12").Formula = "=SUM(D5:D10)"

# --- 6. Number formatting: 0.00% on the new Percentage cells ----------------
$ws.Range("D5:D10").NumberFormat = "0.00%"
$ws.Range("D12").NumberFormat = "0.00%"

# --- 7. Column D width (matches the "Percentage" header column) -------------
$ws.Columns("D").ColumnWidth = 11

# --- 8. Clear the black decorative fill in column E (rows 4-19) so it reads
# as a plain cell again, while columns F:M keep the original black fill -----
$ws.Range("E4:E19").Interior.Pattern = -4142  ## xlPatternNone

# --- 9. Selection / active cell, like the author left it --------------------
$ws.Range("B15").Select()

# --- 10. Point the pie chart at the relocated table ---------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.SeriesCollection(1).Formula = "=SERIES(Sheet1!`$C`$4,Sheet1!`$B`$5:`$B`$10,Sheet1!`$C`$5:`$C`$10,1)"

# --- 11. Reposition/resize the chart frame to its new anchor -----------------
# (from col F/row4 offset to col M/row18 offset, in points, matching the
# sheet's default column width / row height model)
$co.Left = 437.74700787401576
$co.Top = 56.62496062992126
$co.Width = 388.8123622047244
$co.Height = 212.62503937007875

# --- 12. Page orientation -----------------------------------------------------
$ws.PageSetup.Orientation = 1   ## xlPortrait
